$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing trailing zeros / exact text).
$textForceRows = @(5, 6, 10, 11, 16, 21, 23, 24, 26, 28, 30, 31, 32, 33, 34, 36, 37, 38, 40, 41, 46, 47, 48, 49)
foreach ($r in $textForceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.652.69'
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('D3').Value = '2.283.24'
$ws.Range('E3').Value = '  -3.76%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '300.96'
$ws.Range('E5').Value = '  -2.98%  '
$ws.Range('D6').Value = '97.45'
$ws.Range('E6').Value = '  -6.36%  '
$ws.Range('E7').Value = '  -1.76%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -3.93%  '
$ws.Range('D10').Value = '33.55'
$ws.Range('E10').Value = '  -6.07%  '
$ws.Range('D11').Value = '50.78'
$ws.Range('E11').Value = '  -4.63%  '
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('E14').Value = '  -4.04%  '
$ws.Range('D15').Value = '2.636.48'
$ws.Range('E15').Value = '  -3.85%  '
$ws.Range('D16').Value = '15.29'
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('D17').Value = '2.281.22'
$ws.Range('E17').Value = '  -3.85%  '
$ws.Range('E18').Value = '  -2.53%  '
$ws.Range('D19').Value = '42.531.35'
$ws.Range('E19').Value = '  -2.14%  '
$ws.Range('D20').Value = '0.0₃0895'
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('D21').Value = '11.47'
$ws.Range('E21').Value = '  -3.51%  '
$ws.Range('E22').Value = '  -5.03%  '
$ws.Range('D23').Value = '66.68'
$ws.Range('D24').Value = '235.21'
$ws.Range('E24').Value = '  -2.16%  '
$ws.Range('E25').Value = '  -5.00%  '
$ws.Range('D26').Value = '2.48'
$ws.Range('E26').Value = '  -4.67%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = '24.51'
$ws.Range('E28').Value = '  -4.91%  '
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('D30').Value = '164.76'
$ws.Range('E30').Value = '  +2.02%  '
$ws.Range('D31').Value = '33.64'
$ws.Range('E31').Value = '  -7.89%  '
$ws.Range('D32').Value = '9.10'
$ws.Range('E32').Value = '  -3.87%  '
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '4.97'
$ws.Range('E34').Value = '  -4.71%  '
$ws.Range('E35').Value = '  -4.42%  '
$ws.Range('D36').Value = '0.0695'
$ws.Range('E36').Value = '  -5.46%  '
$ws.Range('D37').Value = '4.35'
$ws.Range('E37').Value = '  -6.66%  '
$ws.Range('D38').Value = '16.25'
$ws.Range('E38').Value = '  -10.66%  '
$ws.Range('E39').Value = '  -8.20%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').Value = '1.77'
$ws.Range('E40').Value = '  -7.94%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.0997'
$ws.Range('E41').Value = '  -5.26%  '
$ws.Range('E42').Value = '  -3.15%  '
$ws.Range('E43').Value = '  -7.96%  '
$ws.Range('D44').Value = '1.960.00'
$ws.Range('E44').Value = '  -3.59%  '
$ws.Range('E45').Value = '  -2.45%  '
$ws.Range('D46').Value = '17.72'
$ws.Range('E46').Value = '  -9.48%  '
$ws.Range('D47').Value = '9.68'
$ws.Range('E47').Value = '  -8.11%  '
$ws.Range('D48').Value = '2.84'
$ws.Range('E48').Value = '  -8.43%  '
$ws.Range('D49').Value = '53.36'
$ws.Range('E49').Value = '  -7.39%  '
$ws.Range('E50').Value = '  -3.71%  '
$ws.Range('D51').Value = '2.504.68'
$ws.Range('E51').Value = '  -4.08%  '
